$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$ws.Range("A1").Value = "Datos actualizados a 2 de Abril de 2020 a las 14:50"
$ws.Range("B8").Value = 79465
$ws.Range("C8").Value = 1484
$ws.Range("E8").Value = 59331
$ws.Range("G8").Value = 28
$ws.Range("H8").Value = 959
$ws.Range("B16").Value = 10927
$ws.Range("C16").Value = 216
$ws.Range("E16").Value = 9020
$ws.Range("B21").Value = 6360
$ws.Range("C21").Value = 268
$ws.Range("E21").Value = 6038
$ws.Range("D41").Value = 328
$ws.Range("E41").Value = 1371
$ws.Range("G41").Value = 5
$ws.Range("H41").Value = 21
$ws.Range("E42").Value = 1199
$ws.Range("G42").Value = 2
$ws.Range("H42").Value = 19
$ws.Range("E43").Value = 1311
$ws.Range("G43").Value = 1
$ws.Range("H43").Value = 52
$ws.Range("A87").Value = "Burkina Faso"
$ws.Range("B87").Value = 288
$ws.Range("C87").Value = 6
$ws.Range("D87").Value = 50
$ws.Range("E87").Value = 222
$ws.Range("G87").Value = 0
$ws.Range("H87").Value = 16
$ws.Range("A88").Value = "Camerun"
$ws.Range("B88").Value = 284
$ws.Range("C88").Value = 51
$ws.Range("D88").Value = 10
$ws.Range("E88").Value = 267
$ws.Range("G88").Value = 1
$ws.Range("H88").Value = 7
$ws.Range("D102").Value = 25
$ws.Range("E102").Value = 163
$ws.Range("A119").Value = "Kenia"
$ws.Range("C119").Value = 29
$ws.Range("D119").Value = 4
$ws.Range("E119").Value = 103
$ws.Range("F119").Value = 2
$ws.Range("G119").Value = 2
$ws.Range("H119").Value = 3
$ws.Range("A120").Value = "Camboya"
$ws.Range("B120").Value = 110
$ws.Range("C120").Value = 1
$ws.Range("D120").Value = 34
$ws.Range("E120").Value = 76
$ws.Range("F120").Value = 1
$ws.Range("H120").Value = 0
$ws.Range("A121").Value = "Trinidad yTobago"
$ws.Range("B121").Value = 90
$ws.Range("D121").Value = 1
$ws.Range("E121").Value = 84
$ws.Range("H121").Value = 5
$ws.Range("A122").Value = "Ruanda"
$ws.Range("B122").Value = 82
$ws.Range("D122").Value = 0
$ws.Range("E122").Value = 82
$ws.Range("F122").Value = 0
$ws.Range("H122").Value = 0
$ws.Range("B128").Value = 59
$ws.Range("C128").Value = 2
$ws.Range("E128").Value = 59
$ws.Range("A139").Value = "Republica de Yibuti"
$ws.Range("B139").Value = 40
$ws.Range("C139").Value = 7
$ws.Range("E139").Value = 40
$ws.Range("G139").Value = 0
$ws.Range("H139").Value = 0
$ws.Range("A140").Value = "Zambia"
$ws.Range("C140").Value = 3
$ws.Range("D140").Value = 0
$ws.Range("E140").Value = 38
$ws.Range("G140").Value = 1
$ws.Range("H140").Value = 1
$ws.Range("A141").Value = "Puerto Rico"
$ws.Range("B141").Value = 39
$ws.Range("D141").Value = 1
$ws.Range("E141").Value = 36
$ws.Range("F141").Value = 0
$ws.Range("H141").Value = 2
$ws.Range("A142").Value = "Polinesia Francesa"
$ws.Range("B142").Value = 37
$ws.Range("D142").Value = 0
$ws.Range("E142").Value = 37
$ws.Range("F142").Value = 1
$ws.Range("H142").Value = 0
$ws.Range("A143").Value = "Togo"
$ws.Range("B143").Value = 36
$ws.Range("D143").Value = 10
$ws.Range("E143").Value = 24
$ws.Range("H143").Value = 2
$ws.Range("D153").Value = 2
$ws.Range("E153").Value = 17
